$d = $word.ActiveDocument

# Locate the placeholder run that follows "UBICACIÓN FINAL:" — i.e. the
# " {{Ubicación}}" text.
$target = $d.Content
$found = $target.Find.Execute(" {{Ubicación}}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the ' {{Ubicación}}' placeholder run."
}

# Walk the document's paragraph collection to find the paragraph that owns
# the match, so we can read back its own text (the untouched
# "UBICACIÓN FINAL:" prefix) and its identity attributes (paraId / rsids /
# …) straight from the live document rather than assuming fixed values.
$prefixText = ""
$paraAttrs = ""
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($target.Start -ge $p.Range.Start -and $target.Start -lt $p.Range.End) {
        $prefixRange = $d.Range($p.Range.Start, $target.Start)
        $prefixText = $prefixRange.Text

        $paraXml = $p.Range.WordOpenXML
        if ($paraXml -match '<w:p ([^>]*)>') {
            $paraAttrs = " " + $matches[1]
        }
        break
    }
}

# Only flag xml:space="preserve" on the untouched prefix run when it
# actually has leading/trailing whitespace worth preserving — matching how
# this run already looked (no such attribute) before the edit.
if ($prefixText -match '^\s|\s$') {
    $prefixSpacePreserve = ' xml:space="preserve"'
} else {
    $prefixSpacePreserve = ''
}

# Remove the old placeholder run, then splice in the new run / proofErr
# layout: the literal "{{" in its own run, and "Ubicacion" (flagged by the
# spell-checker once the accent is dropped, hence wrapped in
# spellStart/spellEnd proofErr markers) split across three runs — mirroring
# how Word itself fragments text that proofing touches.
$target.Delete()

$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' +
  '<w:p' + $paraAttrs + '>' +
  '<w:r><w:t' + $prefixSpacePreserve + '>' + $prefixText + '</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> {{</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Ubicaci</w:t></w:r>' +
  '<w:r><w:t>o</w:t></w:r>' +
  '<w:r><w:t>n</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>}}</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$null = $target.InsertXML($newParaXml)
